{"js": "// Locate the two target paragraphs by matching their current text, then:\n//  1. Rewrite the \"sand down the top\" bullet's text and demote it to ilvl 1.\n//  2. Demote the final \"pogo pins\" bullet to ilvl 1.\n//  3. Insert a brand-new top-level (ilvl 0) bullet right after it.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nconst sandParaIndex = items.findIndex(p =>\n  p.text.indexOf(\"Probably need to sand down the top\") !== -1\n);\nif (sandParaIndex === -1) {\n  throw new Error(\"Could not find the 'sand down the top' paragraph\");\n}\nconst sandPara = items[sandParaIndex];\nsandPara.insertText(\n  \"Sanding down the top with an abrasive pad to clean up rough parts: don\\u2019t use anything too abrasive because it can scratch the copper\",\n  Word.InsertLocation.replace\n);\nsandPara.listItem.level = 1;\n\nconst pogoParaIndex = items.findIndex(p =>\n  p.text.indexOf(\"Through holes for pogo pins are far too big\") !== -1\n);\nif (pogoParaIndex === -1) {\n  throw new Error(\"Could not find the 'pogo pins' paragraph\");\n}\nconst pogoPara = items[pogoParaIndex];\n\n// Insert the new bullet right after the pogo-pins paragraph *before* we\n// touch its own indent level, so the new paragraph inherits ilvl 0 (the\n// level the pogo-pins paragraph currently has) rather than the demoted one.\nconst newPara = pogoPara.insertParagraph(\n  \"Looks like the h-bridge driver pads run up on each other\",\n  Word.InsertLocation.after\n);\nnewPara.listItem.level = 0;\n\npogoPara.listItem.level = 1;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Edit 1: rewrite the \"sand down the top\" bullet's text and demote it\n#     from ilvl 0 (ListLevelNumber 1) to ilvl 1 (ListLevelNumber 2). ---\n$rng1 = $d.Content\n[void]$rng1.Find.Execute(\"Probably need to sand down the top\")\n$sandPara = $rng1.Paragraphs(1)\n$newText = \"Sanding down the top with an abrasive pad to clean up rough parts: don\" + [char]0x2019 + \"t use anything too abrasive because it can scratch the copper\"\n$sandPara.Range.Text = $newText\n$sandPara.Range.ListFormat.ListLevelNumber = 2\n\n# --- Edit 2 & 3: add a new top-level (ilvl 0) bullet right after the\n#     \"pogo pins\" bullet, then demote the \"pogo pins\" bullet itself to\n#     ilvl 1. The new paragraph is inserted *before* the demote so it\n#     inherits the still-current ilvl 0 formatting, matching the target. ---\n$rng2 = $d.Content\n[void]$rng2.Find.Execute(\"Through holes for pogo pins are far too big\")\n$pogoPara = $rng2.Paragraphs(1)\n$pogoPara.Range.InsertParagraphAfter()\n\n$newCount = $d.Paragraphs.Count\n$newPara = $d.Paragraphs($newCount)\n$newPara.Range.Text = \"Looks like the h-bridge driver pads run up on each other\"\n\n$pogoPara.Range.ListFormat.ListLevelNumber = 2\n"}
